$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 43

$ws.Cells.Item($row, 1).Value = 53
# Date column stores plain text like "2026-02-16" (not a real Excel date) in
# this sheet, so force text interpretation with a leading quote, same as
# typing '2026-02-16 into the cell in Excel.
$ws.Cells.Item($row, 2).Value = "'2026-02-16"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Value = "21:31:17"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"
$ws.Cells.Item($row, 6).Value = 68702.77
$ws.Cells.Item($row, 7).Value = ""
$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.75
$ws.Cells.Item($row, 12).Value = "Coinbase leading with -0.078% move"
$ws.Cells.Item($row, 13).Value = ""
$ws.Cells.Item($row, 14).Value = 0
